# Refresh the cryptocurrency price / 1h-volume columns (D, E) with the
# latest scraped figures. D-column values that look like plain numbers
# (e.g. "329.36") are written with a leading single-quote so Excel keeps
# them as literal text (matching the original "inline string" cell type)
# instead of silently converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.084.68'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '2.391.07'
$ws.Range("E3").Value = '  +6.35%  '
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").Value = '''329.36'
$ws.Range("E5").Value = '  +10.83%  '
$ws.Range("D6").Value = '''104.13'
$ws.Range("E6").Value = '  -7.51%  '
$ws.Range("E7").Value = '  +3.04%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '''0.649'
$ws.Range("E9").Value = '  +6.87%  '
$ws.Range("D10").Value = '''41.85'
$ws.Range("E10").Value = '  -4.87%  '
$ws.Range("D11").Value = '''0.0937'
$ws.Range("E11").Value = '  +1.18%  '
$ws.Range("E12").Value = '  -3.49%  '
$ws.Range("E13").Value = '  -2.61%  '
$ws.Range("D14").Value = '''17.17'
$ws.Range("E14").Value = '  +12.64%  '
$ws.Range("E15").Value = '  +1.93%  '
$ws.Range("D16").Value = '2.749.81'
$ws.Range("E16").Value = '  +6.27%  '
$ws.Range("D17").Value = '2.390.54'
$ws.Range("E17").Value = '  +4.96%  '
$ws.Range("D18").Value = '43.076.93'
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").Value = '''7.76'
$ws.Range("E19").Value = '  +8.03%  '
$ws.Range("E20").Value = '  +2.04%  '
$ws.Range("D21").Value = '''76.48'
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("E22").Value = '  +4.61%  '
$ws.Range("D23").Value = '''271.68'
$ws.Range("E23").Value = '  +7.72%  '
$ws.Range("E24").Value = '  -1.45%  '
$ws.Range("D25").Value = '''9.68'
$ws.Range("E25").Value = '  +8.09%  '
$ws.Range("D26").Value = '''11.75'
$ws.Range("E26").Value = '  +1.56%  '
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("D28").Value = '''22.97'
$ws.Range("E28").Value = '  +3.65%  '
$ws.Range("E29").Value = '  -1.97%  '
$ws.Range("D30").Value = '''175.17'
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("D31").Value = '''37.49'
$ws.Range("E31").Value = '  -1.26%  '
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").Value = '''0.0929'
$ws.Range("E33").Value = '  +4.37%  '
$ws.Range("D34").Value = '''5.91'
$ws.Range("E34").Value = '  +4.23%  '
$ws.Range("E35").Value = '  +5.10%  '
$ws.Range("D36").Value = '''4.94'
$ws.Range("E36").Value = '  -2.63%  '
$ws.Range("D37").Value = '''4.15'
$ws.Range("E37").Value = '  -2.04%  '
$ws.Range("E38").Value = '  -2.54%  '
$ws.Range("D39").Value = '''0.108'
$ws.Range("E39").Value = '  +3.28%  '
$ws.Range("D40").Value = '''2.79'
$ws.Range("E40").Value = '  +15.79%  '
$ws.Range("D41").Value = '''1.58'
$ws.Range("E41").Value = '  +19.53%  '
$ws.Range("E42").Value = '  +0.88%  '
$ws.Range("D43").Value = '''69.71'
$ws.Range("E43").Value = '  -3.34%  '
$ws.Range("D44").Value = '''121.54'
$ws.Range("E44").Value = '  +14.16%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").Value = '''12.38'
$ws.Range("E46").Value = '  -0.87%  '
$ws.Range("D47").Value = '''90.44'
$ws.Range("E47").Value = '  +47.46%  '
$ws.Range("D48").Value = '''9.35'
$ws.Range("E48").Value = '  +7.98%  '
$ws.Range("E49").Value = '  +0.30%  '
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("D51").Value = '''0.496'
$ws.Range("E51").Value = '  +12.84%  '
